$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New id (#slug) / speaker_variant pairs per row, re-exported with
# the id freshly derived from speaker_variant (no levenshtein-matched id)
# and with is_prefered left blank for every row.
$rows = @(
    @{ Row = 2;  Id = "#roxane";       Variant = "Roxane" },
    @{ Row = 3;  Id = "#ridderen";     Variant = "Ridderen" },
    @{ Row = 4;  Id = "#boode";        Variant = "Boode" },
    @{ Row = 5;  Id = "#arbates";      Variant = "Arbates" },
    @{ Row = 6;  Id = "#antithakata";  Variant = "Antithakata" },
    @{ Row = 7;  Id = "#krater";       Variant = "Krater" },
    @{ Row = 8;  Id = "#hofknaap";     Variant = "Hofknaap" },
    @{ Row = 9;  Id = "#eumelus";      Variant = "Eumelus" },
    @{ Row = 10; Id = "#arbate";       Variant = "Arbate" },
    @{ Row = 11; Id = "#admetus";      Variant = "Admetus" },
    @{ Row = 12; Id = "#philippus";    Variant = "Philippus" },
    @{ Row = 13; Id = "#antithanata";  Variant = "Antithanata" },
    @{ Row = 14; Id = "#admeius";      Variant = "Admeius" },
    @{ Row = 15; Id = "#gevolg-van";   Variant = "Gevolg van" },
    @{ Row = 16; Id = "#alcestis";     Variant = "Alcestis" },
    @{ Row = 17; Id = "#anthithanata"; Variant = "Anthithanata" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.Id
    $ws.Cells.Item($r.Row, 3).Value = $r.Variant
    $ws.Cells.Item($r.Row, 4).ClearContents()
}
